# Remove the surname " Orsi" from "contribuição do Fábio Orsi." leaving
# "contribuição do Fábio." — while keeping the final "o" of "Fábio" in its
# own run (matching the target canonical OOXML run split: "...Fábi" / "o").

$d = $word.ActiveDocument

# Locate "Fábio Orsi" (unique in the document) and remember its extent.
$hit = $d.Content
$found = $hit.Find.Execute("Fábio Orsi", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if ($found) {
    $matchStart = $hit.Start
    $matchEnd   = $hit.End

    $prefixLen = "Fábi".Length          # chars that stay in the first run
    $oStart    = $matchStart + $prefixLen
    $oEnd      = $oStart + 1            # just the trailing "o" of "Fábio"

    # Toggle a character property on just the "o" so Word is forced to split
    # it into its own run, isolated from the surrounding text, instead of
    # merging it back into the identically-formatted neighbouring runs.
    $oRange = $d.Range($oStart, $oEnd)
    $oRange.Bold = 1

    # Delete everything from right after that "o" through the end of the
    # match (i.e. " Orsi"), while the "o" is still set apart as its own run.
    $tailRange = $d.Range($oEnd, $matchEnd)
    $tailRange.Text = ""

    # Restore the "o" run's formatting back to normal (not bold) — it keeps
    # living as its own run, just like the one " Orsi" used to occupy.
    $oRange2 = $d.Range($oStart, $oStart + 1)
    $oRange2.Bold = 0
}
